$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B13 ("Implement merge nodes") now marked as done -> apply "Good" cell style
$ws.Range("B13").Style = "Good"

# --- "hours spent" block rework (rows 26-37, total moved from row 38 to row 42) ---

# Row 26: the "hours spent" label moves up from A27 to A26
$ws.Range("A26").Value = "hours spent"

# Row 27
$ws.Range("A27").Value = 2
$ws.Range("B27").Value = 36

# Row 28
$ws.Range("A28").Value = 3
$ws.Range("B28").Value = 36

# Row 29
$ws.Range("A29").Value = 4
$ws.Range("B29").Value = 20

# Row 30
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = 31

# Row 31
$ws.Range("A31").Value = 6
$ws.Range("B31").Value = 20

# Row 32
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = 26

# Row 33
$ws.Range("A33").Value = 8
$ws.Range("B33").Value = 20

# Row 34
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = 22

# Row 35 (new)
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = 14

# Row 36 (new)
$ws.Range("A36").Value = 13
$ws.Range("B36").Value = 30

# Row 37 (new)
$ws.Range("A37").Value = 14
$ws.Range("B37").Formula = "=14+19"

# Clear out the old totals row (row 38 no longer used)
$ws.Range("A38").Value = ""
$ws.Range("B38").Value = ""
$ws.Range("C38").Value = ""

# New totals row 42
$ws.Range("A42").Value = "Total hours spent so far"
$ws.Range("B42").Formula = "=SUM(B27:B40)"
$ws.Range("C42").Formula = "=B42*135"
